# Project DesignFirst is saved. Author: admin. Type: SAVE.
# Semantic change: cell C10 on the "Rules" sheet changes from 18 to 100.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Rules")

$ws.Range("C10").Value = 100
